# OncoTree, DO, EFO updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# source_version column (E) updates for OncoTree, Disease Ontology, EFO
$ws.Range("E2").Value = "2025_10_03"   # OncoTree version
$ws.Range("E3").Value = "v2025-11-25"  # Disease Ontology version
$ws.Range("E4").Value = "v3.84.0"      # EFO version

# Move selection to E2 (also clears the scrolled topLeftCell state)
$ws.Range("E2").Select()
